# Fine tune the intents, save report function
# Rename headers, drop unused columns (router_status, phone_status, Instrument_status)
# and keep the final working status in column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing columns (D:F) entirely, shifting nothing else.
$ws.Range("D1:F2").EntireColumn.Delete()

# Update header labels
$ws.Range("B1").Value = "phoneNumber"
$ws.Range("C1").Value = "issue"

# Update data row: the remaining status value
$ws.Range("C2").Value = "Working"
